$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "IncorrectUser"
$ws.Range("A2").Value = "incorrect uname"
